# Requirements workbook update:
#  - split the sheet into "requirements" and a new "nonfunctional" sheet
#  - move 4 non-functional-ish requirements (rq_dup_req_search, rq_dup_same_version,
#    rq_markup_remove, rq_automatic_diagram) to the new sheet
#  - tweak the rq_automatic_diagram description text
#  - fix up row heights, autofilter range and the hidden _FilterDatabase names

$wb = $excel.ActiveWorkbook
$reqs = $wb.Worksheets.Item(1)
$reqs.Name = "requirements"

# New sheet, placed right after "requirements"
$nonf = $wb.Worksheets.Add($null, $reqs)
$nonf.Name = "nonfunctional"

# Header row (A1:G1) copied as-is onto the new sheet
$reqs.Range("A1:G1").Copy()
$nonf.Range("A1").PasteSpecial()
$nonf.Rows.Item(1).RowHeight = 13

# The four rows being relocated (old row -> new row on "nonfunctional")
$moveRows = @(73, 74, 75, 76)
$destRow = 2
foreach ($srcRow in $moveRows) {
    $reqs.Range("A" + $srcRow + ":E" + $srcRow).Copy()
    $nonf.Range("A" + $destRow).PasteSpecial()
    $destRow = $destRow + 1
}

# Reword the rq_automatic_diagram description (now row 5 on "nonfunctional")
$nonf.Range("D5").Value = "VRM2 shall be able to generate diagrams (requirements, hierarchy, safety) from the command line and with all relevant pameters specified on the command line." + [char]10 + "Comment: This feature is intended for CI where automatic generation of a diagram is desirable."

# New row heights on "nonfunctional"
$nonf.Rows.Item(2).RowHeight = 90
$nonf.Rows.Item(3).RowHeight = 57
$nonf.Rows.Item(4).RowHeight = 57
$nonf.Rows.Item(5).RowHeight = 101

# Remove the relocated rows from "requirements" (delete bottom-up so row numbers stay valid)
$reqs.Rows.Item(76).Delete()
$reqs.Rows.Item(75).Delete()
$reqs.Rows.Item(74).Delete()
$reqs.Rows.Item(73).Delete()

# Row heights that shifted on "requirements" because the trailing rows moved away
$reqs.Rows.Item(53).RowHeight = 24
$reqs.Rows.Item(55).RowHeight = 46
$reqs.Rows.Item(56).RowHeight = 24
$reqs.Rows.Item(57).RowHeight = 13
$reqs.Rows.Item(59).RowHeight = 13
$reqs.Rows.Item(60).RowHeight = 13
$reqs.Rows.Item(66).RowHeight = 24
$reqs.Rows.Item(70).RowHeight = 13
$reqs.Rows.Item(71).RowHeight = 24
$reqs.Rows.Item(72).RowHeight = 35

# Autofilter now spans only the remaining 72 rows
$reqs.Range("A1:G72").AutoFilter()

# Hidden filter-database defined names need to point at the renamed sheet / new extent
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=requirements!`$A`$1:`$G`$72"
    }
    if ($n.Name -eq "Sheet1!_FilterDatabase_0") {
        $n.RefersTo = "=requirements!`$A`$1:`$G`$70"
    }
}

# Restore selection on the requirements sheet
$reqs.Range("A54").Select()
$nonf.Range("A2").Select()
$reqs.Select()
